$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")
$lo = $ws.ListObjects.Item(1)

# Insert two new rows right after the header row (shifts old rows 2-14 down to 4-16).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(3).Insert()

# Populate the two new rows with copied localization-key data.
# (Set A2/A3/B3 now; B2 is set later so the shared-string table ends up
# ordered the same way the original authoring tool produced it.)
$ws.Range("A2").Value = "competition-key"
$ws.Range("A3").Value = "host-key"
$ws.Range("B3").Value = "germany"

# Append ten new rows at the bottom with venue-key / host-venue pairs.
$ws.Range("A17").Value = "venue-key.1"
$ws.Range("B17").Value = "de-berlin"
$ws.Range("A18").Value = "venue-key.2"
$ws.Range("B18").Value = "de-cologne"
$ws.Range("A19").Value = "venue-key.3"
$ws.Range("B19").Value = "de-dortmund"
$ws.Range("A20").Value = "venue-key.4"
$ws.Range("B20").Value = "de-dusseldorf"
$ws.Range("A21").Value = "venue-key.5"
$ws.Range("B21").Value = "de-frankfurt"
$ws.Range("A22").Value = "venue-key.6"
$ws.Range("B22").Value = "de-gelsenkirchen"
$ws.Range("A23").Value = "venue-key.7"
$ws.Range("B23").Value = "de-hamburg"
$ws.Range("A24").Value = "venue-key.8"
$ws.Range("B24").Value = "de-leipzig"
$ws.Range("A25").Value = "venue-key.9"
$ws.Range("B25").Value = "de-munich"
$ws.Range("A26").Value = "venue-key.10"
$ws.Range("B26").Value = "de-stuttgart"

# Finally set B2 (added last so it becomes the last new shared string).
$ws.Range("B2").Value = "mens-euro"

# Resize the "tournament" table to cover the new rows.
$lo.Resize($ws.Range("A1:I26"))

# Update the selected cell to match the saved view state.
$ws.Range("B3").Select()
